$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the validation string for row 7 (J7): "notification" -> "notifications"
$ws.Range("J7").Value = "status=200||notifications=true"

# Clear the STATUS column (L) for all data rows - the PASS markers are removed
$ws.Range("L2:L10").ClearContents()

# Clear stray empty placeholder cells in DEPENDENCYTESTS (I) and STORE (K) columns
$ws.Range("I5").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("K10").ClearContents()

# Update selection / view to match
$ws.Application.ActiveWindow.ScrollColumn = 8
$ws.Range("L2:L10").Select()
